# Insert a new weekly data row for "Terminal Hortofrutícola Agro Chillán - Piña"
# at sheet row 388, pushing the existing rows 388:405 down to 389:406.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(388).Insert()

$ws.Cells.Item(388, 1).Value = 7
$ws.Cells.Item(388, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(388, 3).Value = "Ñuble"
$ws.Cells.Item(388, 4).Value = 45267
$ws.Cells.Item(388, 5).Value = 16
$ws.Cells.Item(388, 6).Value = "Fruta"
$ws.Cells.Item(388, 7).Value = 100108
$ws.Cells.Item(388, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(388, 9).Value = 100108005
$ws.Cells.Item(388, 10).Value = "Piña"
$ws.Cells.Item(388, 11).Value = "Caramelo"
$ws.Cells.Item(388, 12).Value = "Segunda"
$ws.Cells.Item(388, 13).Value = 180
$ws.Cells.Item(388, 14).Value = 26000
$ws.Cells.Item(388, 15).Value = 27000
$ws.Cells.Item(388, 16).Value = 26556
$ws.Cells.Item(388, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(388, 18).Value = "Ecuador"
$ws.Cells.Item(388, 19).Value = 1897
$ws.Cells.Item(388, 20).Value = 14
